# Generate Report for Handoff
#
# The localization status report is regenerated: the "Handed back: in sync
# with en-US" status becomes "Ready for handoff" everywhere it is used, and
# the associated timestamps are refreshed. The now-narrower status text also
# lets the status columns be resized down from their old (wider) width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps ---
$wsOverview.Range("G2").Value = "2016-09-05 15:14:08"
$wsDeDe.Range("H2").Value     = "2016-09-05 15:14:08"
$wsZhCn.Range("H2").Value     = "2016-09-05 15:13:58"

# --- Status columns shrink to fit the shorter "Ready for handoff" text ---
# (target stored column width ~17.216 chars; the host quantizes ColumnWidth
# to 1/6-character pixel buckets, so 16.16 is the closest input that lands
# on the nearest reachable bucket)
$wsOverview.Columns.Item(5).ColumnWidth = 16.16
$wsOverview.Columns.Item(6).ColumnWidth = 16.16
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.16
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.16
